$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - TC_006
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "TC_006"
$ws.Range("C9").Value = "standard_user"
$ws.Range("D9").Value = "secret_sauce"

# Row 10 - TC_007
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "TC_007"
$ws.Range("C10").Value = "standard_user"
$ws.Range("D10").Value = "secret_sauce"
$ws.Range("I1").Value = "ExpectedCount"
$ws.Range("I10").Value = 6

# Row 11 - TC_008
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "TC_008"
$ws.Range("C11").Value = "standard_user"
$ws.Range("D11").Value = "secret_sauce"
$ws.Range("J1").Value = "PriceTag"
$ws.Range("J11").Value = "$"

# Row 12 - TC_009
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "TC_009"
$ws.Range("C12").Value = "standard_user"
$ws.Range("D12").Value = "secret_sauce"
$ws.Range("K1").Value = "SortValue"
$ws.Range("K12").Value = "az"

# Match the resulting selection state left behind in the workbook
$ws.Range("K13").Select()
